{"js": "// Remove the phone number from the contact-details line near the top of\n// the CV: \"Worcestershire | +447494856994 | richards.kimc@gmail.com\"\n// becomes \"Worcestershire | richards.kimc@gmail.com\".\nconst body = context.document.body;\n\nconst results = body.search(\"+447494856994 | \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const range of results.items) {\n  range.insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Remove the phone number from the contact-details line near the top of\n# the CV: \"Worcestershire | +447494856994 | richards.kimc@gmail.com\"\n# becomes \"Worcestershire | richards.kimc@gmail.com\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"+447494856994 | \"\n$find.Replacement.Text = \"\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
